$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 61820.59
$ws.Range("I64").Value = 127625
$ws.Range("K64").Value = 127625
$ws.Range("M64").Value = -127377

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 61820.59
$ws.Range("I67").Value = 127625
$ws.Range("K67").Value = 127625
$ws.Range("M67").Value = -126767

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 50001156
$ws.Range("J127").Value = 90910740
$ws.Range("L127").Value = 272732220
$ws.Range("N127").Value = -272742140

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3950.4
$ws.Range("I138").Value = 3756.6365
$ws.Range("J138").Value = 3981.2898
$ws.Range("K138").Value = 11269.9095
$ws.Range("L138").Value = 11943.8694
$ws.Range("M138").Value = -6129.9095
$ws.Range("N138").Value = -22223.8694

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 50375.52
$ws.Range("I32").Value = 20830.158
$ws.Range("J32").Value = 212875
$ws.Range("K32").Value = 20830.158
$ws.Range("L32").Value = 212875
$ws.Range("M32").Value = -20543.158
$ws.Range("N32").Value = -213449

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1867.25
$ws.Range("I74").Value = 1910.6
$ws.Range("J74").Value = 1795
$ws.Range("K74").Value = 1910.6
$ws.Range("L74").Value = 1795
$ws.Range("M74").Value = -1036.6
$ws.Range("N74").Value = -3543

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1867.25
$ws.Range("I77").Value = 1910.6
$ws.Range("J77").Value = 1795
$ws.Range("K77").Value = 9553
$ws.Range("L77").Value = 8975
$ws.Range("M77").Value = -5185
$ws.Range("N77").Value = -17711

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1665.7273
$ws.Range("I88").Value = 1217
$ws.Range("J88").Value = 1765.4445
$ws.Range("K88").Value = 1217
$ws.Range("L88").Value = 1765.4445
$ws.Range("M88").Value = -811
$ws.Range("N88").Value = -2577.4445

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1665.7273
$ws.Range("I91").Value = 1217
$ws.Range("J91").Value = 1765.4445
$ws.Range("K91").Value = 1217
$ws.Range("L91").Value = 1765.4445
$ws.Range("M91").Value = 187
$ws.Range("N91").Value = -4573.4445

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 44884.26
$ws.Range("I97").Value = 48824.57
$ws.Range("K97").Value = 48824.57
$ws.Range("M97").Value = -48328.57

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10218298
$ws.Range("I132").Value = 11643577
$ws.Range("J132").Value = 3799.1667
$ws.Range("K132").Value = 34930731
$ws.Range("L132").Value = 11397.5001
$ws.Range("M132").Value = -34928201
$ws.Range("N132").Value = -16457.5001

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 21168.857
$ws.Range("J35").Value = 21168.857
$ws.Range("L35").Value = 21168.857
$ws.Range("N35").Value = -21788.857

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17864.572
$ws.Range("J82").Value = 28890
$ws.Range("L82").Value = 28890
$ws.Range("N82").Value = -29656

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 17864.572
$ws.Range("J85").Value = 28890
$ws.Range("L85").Value = 28890
$ws.Range("N85").Value = -31542

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 87765.46000000001
$ws.Range("I86").Value = 113629.4
$ws.Range("K86").Value = 113629.4
$ws.Range("M86").Value = -112506.4

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 87765.46000000001
$ws.Range("I89").Value = 113629.4
$ws.Range("K89").Value = 568147
$ws.Range("M89").Value = -562531

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 88847.05
$ws.Range("I105").Value = 78900.69500000001
$ws.Range("J105").Value = 101777.3
$ws.Range("K105").Value = 78900.69500000001
$ws.Range("L105").Value = 101777.3
$ws.Range("M105").Value = -77153.69500000001
$ws.Range("N105").Value = -105271.3

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1145.3334
$ws.Range("I122").Value = 1049.75
$ws.Range("J122").Value = 1221.8
$ws.Range("K122").Value = 3149.25
$ws.Range("L122").Value = 3665.4
$ws.Range("M122").Value = -699.25
$ws.Range("N122").Value = -8565.4

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 700
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 900
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 2700
$ws.Range("M32").Value = -1217
$ws.Range("N32").Value = -3266

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1956.9231
$ws.Range("I81").Value = 931.4286
$ws.Range("J81").Value = 2334.7368
$ws.Range("K81").Value = 2794.2858
$ws.Range("L81").Value = 7004.2104
$ws.Range("M81").Value = -1671.2858
$ws.Range("N81").Value = -9250.2104

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1956.9231
$ws.Range("I84").Value = 931.4286
$ws.Range("J84").Value = 2334.7368
$ws.Range("K84").Value = 8382.857399999999
$ws.Range("L84").Value = 21012.6312
$ws.Range("M84").Value = -2766.857399999999
$ws.Range("N84").Value = -32244.6312

# CUL row 105
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 14614.286
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 16716.666
$ws.Range("K105").Value = 6000
$ws.Range("L105").Value = 50149.99800000001
$ws.Range("M105").Value = -3379
$ws.Range("N105").Value = -55391.99800000001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 939.8333
$ws.Range("I122").Value = 540
$ws.Range("J122").Value = 1499.6
$ws.Range("K122").Value = 4860
$ws.Range("L122").Value = 13496.4
$ws.Range("M122").Value = -2410
$ws.Range("N122").Value = -18396.4

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 805438.0600000001
$ws.Range("J131").Value = 845671.6
$ws.Range("L131").Value = 2537014.8
$ws.Range("N131").Value = -2547094.8

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 29412784
$ws.Range("I97").Value = 35715410
$ws.Range("K97").Value = 35715410
$ws.Range("M97").Value = -35714914

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5231.9375
$ws.Range("I132").Value = 4542.857
$ws.Range("J132").Value = 5767.8887
$ws.Range("K132").Value = 13628.571
$ws.Range("L132").Value = 17303.6661
$ws.Range("M132").Value = -11098.571
$ws.Range("N132").Value = -22363.6661

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 69210.60000000001
$ws.Range("I40").Value = 501999.5
$ws.Range("J40").Value = 2627.6924
$ws.Range("K40").Value = 501999.5
$ws.Range("L40").Value = 2627.6924
$ws.Range("M40").Value = -501863.5
$ws.Range("N40").Value = -2899.6924

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5177.857
$ws.Range("I122").Value = 4588.5293
$ws.Range("K122").Value = 13765.5879
$ws.Range("M122").Value = -11315.5879

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1383.8
$ws.Range("I122").Value = 1484.8
$ws.Range("J122").Value = 979.8
$ws.Range("K122").Value = 4454.4
$ws.Range("L122").Value = 2939.4
$ws.Range("M122").Value = -2004.4
$ws.Range("N122").Value = -7839.4

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1000.6429
$ws.Range("I126").Value = 981.4
$ws.Range("K126").Value = 2944.2
$ws.Range("M126").Value = -474.1999999999998

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19242.186
$ws.Range("I132").Value = 2284.5
$ws.Range("J132").Value = 79246.30499999999
$ws.Range("K132").Value = 6853.5
$ws.Range("L132").Value = 237738.915
$ws.Range("M132").Value = -4323.5
$ws.Range("N132").Value = -242798.915
